# Updated remote access methods to use server_groups + units tests update
#
# Summary of changes applied to Sheet1:
#  - Drop column I (COMPLIANCE_GROUPS) entirely.
#  - Rename header H1 "GROUPS" -> "SERVER_GROUPS".
#  - Row 2 (ssh host entry): hostname -> IP address, creds -> "cyberwatch"/"cyberwatch".
#  - Row 3: unchanged content, but row height / column widths are refreshed.
#  - Row 4: swap the WinRm example row for a second Ssh example row
#           (server01.example.com / 22 / Ssh::WithPassword / master / admin / SuperPassword / production, test).
#  - Column widths tuned (A, C, F widened) and selection moved to row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the COMPLIANCE_GROUPS column (column I) completely.
# ------------------------------------------------------------------
$ws.Columns.Item(9).Delete()

# ------------------------------------------------------------------
# 2. Header row.
# ------------------------------------------------------------------
$ws.Range("H1").Value = "SERVER_GROUPS"

# ------------------------------------------------------------------
# 3. Row 2 - first SSH entry now points at an IP and uses the
#    "cyberwatch" credential pair.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "10.0.2.15"
$ws.Range("F2").Value = "cyberwatch"
$ws.Range("G2").Value = "cyberwatch"

# Row 3 content is unchanged (server02.example.com / ssh key / user).

# ------------------------------------------------------------------
# 4. Row 4 - replace the WinRm example with a second SSH example.
# ------------------------------------------------------------------
$ws.Range("A4").Value = "server01.example.com"
$ws.Range("B4").Value = 22
$ws.Range("C4").Value = "CbwRam::RemoteAccess::Ssh::WithPassword"
$ws.Range("D4").Value = "master"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "admin"
$ws.Range("G4").Value = "SuperPassword"
$ws.Range("H4").Value = "production, test"

# ------------------------------------------------------------------
# 5. Column widths (character units, closest achievable via the
#    ColumnWidth pixel-rounded COM property).
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.75   # A  ~= 20.72
$ws.Columns.Item(3).ColumnWidth = 39.25   # C  ~= 40.24
$ws.Columns.Item(6).ColumnWidth = 20.09   # F  ~= 21.04

# ------------------------------------------------------------------
# 6. Row heights.
# ------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 1572.35

# ------------------------------------------------------------------
# 7. Selection / viewport -> row 4.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(4).Select()
